$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (46061 -> 46062) for every data row, from row 2 through row 320.
for ($r = 2; $r -le 320; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 46061) {
        $cell.Value = 46062
    }
}
